# TC25_Canine_Filter_Breed-GermanShphd.xlsx - "Fixed ICDC breed all testcases"
#
# Replaces the StatQuery Cypher text (column C, rows 2-4) with the
# corrected/shortened query, then updates the view state (zoom, selection,
# row heights) that Excel re-derives when the cell content shrinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['German Shepherd Dog']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# Column C (StatQuery) holds the same query text on rows 2, 3 and 4.
$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery

# Row heights shrink to fit the much shorter query text.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 216

# Window/view state: zoom 70% -> 115%, scrolled to top-left (B1), and the
# selection now spans B4:B5.
$win = $excel.ActiveWindow
$win.Zoom = 115
$win.ScrollRow = 1
$win.ScrollColumn = 2

$ws.Range("B4:B5").Select()
